$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.05
$ws.Range("G2").Value = 3.75
$ws.Range("H2").Value = 2.28
$ws.Range("I2").Value = 2.64
$ws.Range("K2").Value = 3.8
$ws.Range("N2").Value = 3
$ws.Range("P2").Value = 1.73
$ws.Range("Q2").Value = 2.02
$ws.Range("U2").Value = 1.97
$ws.Range("V2").Value = 1.61
$ws.Range("AE2").Value = 80
$ws.Range("AH2").Value = 20
$ws.Range("AO2").Value = 26
$ws.Range("F3").Value = 1.44
$ws.Range("G3").Value = 1.57
$ws.Range("H3").Value = 6.6
$ws.Range("I3").Value = 10
$ws.Range("J3").Value = 3.8
$ws.Range("K3").Value = 5.3
$ws.Range("L3").Value = 1.38
$ws.Range("O3").Value = 1.27
$ws.Range("P3").Value = 1.91
$ws.Range("R3").Value = 1.34
$ws.Range("T3").Value = 1.97
$ws.Range("U3").Value = 1.8
$ws.Range("V3").Value = 1.11
$ws.Range("W3").Value = 2.74
$ws.Range("X3").Value = 1000
$ws.Range("AB3").Value = 16
$ws.Range("AC3").Value = 42
$ws.Range("AF3").Value = 17
$ws.Range("AG3").Value = 28
$ws.Range("AN3").Value = 24
$ws.Range("H4").Value = 8.4
$ws.Range("N4").Value = 3.45
$ws.Range("H5").Value = 6
$ws.Range("J5").Value = 3.8
$ws.Range("K5").Value = 4.7
$ws.Range("N5").Value = 3.05
$ws.Range("T5").Value = 2.02
$ws.Range("U5").Value = 1.76
$ws.Range("V5").Value = 1.14
$ws.Range("AC5").Value = 14
$ws.Range("AF5").Value = 40
$ws.Range("F6").Value = 2.68
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 2.46
$ws.Range("J6").Value = 3.7
$ws.Range("K6").Value = 4.1
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 4.5
$ws.Range("O6").Value = 1.23
$ws.Range("P6").Value = 2.14
$ws.Range("Q6").Value = 1.7
$ws.Range("R6").Value = 1.46
$ws.Range("S6").Value = 2.82
$ws.Range("T6").Value = 1.71
$ws.Range("U6").Value = 2.34
$ws.Range("W6").Value = 1.51
$ws.Range("X6").Value = 22
$ws.Range("Z6").Value = 19.5
$ws.Range("AA6").Value = 900
$ws.Range("AB6").Value = 15
$ws.Range("AC6").Value = 9.6
$ws.Range("AD6").Value = 12.5
$ws.Range("AE6").Value = 65
$ws.Range("AF6").Value = 21
$ws.Range("AG6").Value = 13.5
$ws.Range("AH6").Value = 16
$ws.Range("AI6").Value = 95
$ws.Range("AK6").Value = 75
$ws.Range("AL6").Value = 110
$ws.Range("AM6").Value = 580
$ws.Range("AO6").Value = 19
$ws.Range("F7").Value = 9
$ws.Range("H7").Value = 1.19
$ws.Range("J7").Value = 6.4
$ws.Range("K7").Value = 9.4
$ws.Range("M7").Value = 1.03
$ws.Range("P7").Value = 2.82
$ws.Range("T7").Value = 2.04
$ws.Range("U7").Value = 1.72
$ws.Range("V7").Value = 4.9
$ws.Range("Y7").Value = 13
$ws.Range("F8").Value = 2.7
$ws.Range("H8").Value = 2.58
$ws.Range("I8").Value = 3
$ws.Range("V8").Value = 1.5
$ws.Range("W8").Value = 1.47
$ws.Range("Z8").Value = 19
$ws.Range("AD8").Value = 13.5
$ws.Range("AF8").Value = 21
$ws.Range("AI8").Value = 1000
$ws.Range("AK8").Value = 40
$ws.Range("AN8").Value = 40
$ws.Range("AO8").Value = 600
$ws.Range("H9").Value = 4.8
$ws.Range("N9").Value = 3.75
$ws.Range("S9").Value = 3.5
$ws.Range("T9").Value = 1.86
$ws.Range("U9").Value = 2
$ws.Range("Z9").Value = 40
$ws.Range("AB9").Value = 8.6
$ws.Range("AD9").Value = 20
$ws.Range("AI9").Value = 75
$ws.Range("AK9").Value = 20
$ws.Range("AL9").Value = 36
$ws.Range("AM9").Value = 580
$ws.Range("AO9").Value = 90
$ws.Range("H10").Value = 2.9
$ws.Range("I10").Value = 3.05
$ws.Range("J10").Value = 3.65
$ws.Range("N10").Value = 4.9
$ws.Range("P10").Value = 2.28
$ws.Range("Q10").Value = 1.71
$ws.Range("R10").Value = 1.51
$ws.Range("T10").Value = 1.61
$ws.Range("U10").Value = 2.5
$ws.Range("W10").Value = 1.65
$ws.Range("AA10").Value = 120
$ws.Range("AE10").Value = 32
$ws.Range("H11").Value = 1.89
$ws.Range("S11").Value = 2.36
$ws.Range("U11").Value = 2.6
$ws.Range("X11").Value = 25
$ws.Range("Z11").Value = 15.5
$ws.Range("AB11").Value = 75
$ws.Range("AF11").Value = 36
$ws.Range("AH11").Value = 16
$ws.Range("AI11").Value = 24
$ws.Range("AL11").Value = 120
$ws.Range("AN11").Value = 240
$ws.Range("F12").Value = 2.86
$ws.Range("G12").Value = 3.05
$ws.Range("H12").Value = 2.58
$ws.Range("I12").Value = 2.7
$ws.Range("P12").Value = 1.81
$ws.Range("R12").Value = 1.31
$ws.Range("T12").Value = 1.79
$ws.Range("V12").Value = 1.59
$ws.Range("W12").Value = 1.5
$ws.Range("Y12").Value = 10.5
$ws.Range("Z12").Value = 19.5
$ws.Range("AG12").Value = 15
$ws.Range("AM12").Value = 110
$ws.Range("G13").Value = 1.93
$ws.Range("Q13").Value = 1.66
$ws.Range("F14").Value = 1.9
$ws.Range("G14").Value = 2
$ws.Range("I14").Value = 4.8
$ws.Range("J14").Value = 3.5
$ws.Range("T14").Value = 1.93
$ws.Range("U14").Value = 1.86
$ws.Range("W14").Value = 2
